$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 10) that duplicates row 9's data and formatting,
# mirroring the existing repeated "26-09-2025" entry.
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4104)
